# Data update using gitrun.py
# Updates enrollment ("Inscritos"), payments ("Pagos") and homologated
# registrations ("Inscrições homologadas") counts on the "Inscricoes" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inscricoes")

# Map of cell address -> new value
$updates = @{
    "E2"  = 28
    "F2"  = 14
    "H2"  = 14
    "E9"  = 10
    "E15" = 90
    "F15" = 40
    "H15" = 40
    "F17" = 21
    "H17" = 21
    "E19" = 28
    "E33" = 16
    "E34" = 7
    "E36" = 45
    "F36" = 15
    "H36" = 15
    "E37" = 18
    "E38" = 35
    "E39" = 10
    "E40" = 6
    "E47" = 37
    "F47" = 18
    "H47" = 18
    "E49" = 29
    "E50" = 7
    "E60" = 8
    "E61" = 12
    "E63" = 10
    "E66" = 21
    "E67" = 18
    "E69" = 11
    "E71" = 11
    "F71" = 8
    "H71" = 8
    "E72" = 19
    "E75" = 7
    "E77" = 21
    "E78" = 12
    "E79" = 12
    "E83" = 4
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
